$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 167.14285
$ws.Range("I11").Value = 167.14285
$ws.Range("K11").Value = 167.14285
$ws.Range("M11").Value = -27.14285000000001
$ws.Range("H31").Value = 2626.6667
$ws.Range("I31").Value = 1580.125
$ws.Range("J31").Value = 10999
$ws.Range("K31").Value = 4740.375
$ws.Range("L31").Value = 32997
$ws.Range("M31").Value = -4510.375
$ws.Range("N31").Value = -33457
$ws.Range("H58").Value = 874.5
$ws.Range("J58").Value = 1693.3334
$ws.Range("L58").Value = 5080.0002
$ws.Range("N58").Value = -5380.0002
$ws.Range("I88").Value = 83340080
$ws.Range("J88").Value = 2390698.5
$ws.Range("K88").Value = 83340080
$ws.Range("L88").Value = 2390698.5
$ws.Range("M88").Value = -83339674
$ws.Range("N88").Value = -2391510.5
$ws.Range("I91").Value = 83340080
$ws.Range("J91").Value = 2390698.5
$ws.Range("K91").Value = 83340080
$ws.Range("L91").Value = 2390698.5
$ws.Range("M91").Value = -83338676
$ws.Range("N91").Value = -2393506.5
$ws.Range("H96").Value = 1675.3334
$ws.Range("I96").Value = 2263
$ws.Range("J96").Value = 500
$ws.Range("K96").Value = 6789
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = -5416
$ws.Range("N96").Value = -4246
$ws.Range("H106").Value = 2155.5715
$ws.Range("I106").Value = 2323.1667
$ws.Range("J106").Value = 1150
$ws.Range("K106").Value = 2323.1667
$ws.Range("L106").Value = 1150
$ws.Range("M106").Value = -1692.1667
$ws.Range("N106").Value = -2412
$ws.Range("H116").Value = 6811.2573
$ws.Range("I116").Value = 6770.3667
$ws.Range("K116").Value = 6770.3667
$ws.Range("M116").Value = -3328.3667
$ws.Range("H132").Value = 2875.6785
$ws.Range("I132").Value = 2969.76
$ws.Range("J132").Value = 2091.6667
$ws.Range("K132").Value = 8909.280000000001
$ws.Range("L132").Value = 6275.000100000001
$ws.Range("M132").Value = -6379.280000000001
$ws.Range("N132").Value = -11335.0001
$ws.Range("H137").Value = 2567.487
$ws.Range("I137").Value = 2367.4194
$ws.Range("K137").Value = 7102.2582
$ws.Range("M137").Value = -4552.2582
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1227707.8
$ws.Range("I2").Value = 2453237.2
$ws.Range("J2").Value = 2178.3333
$ws.Range("K2").Value = 2453237.2
$ws.Range("L2").Value = 2178.3333
$ws.Range("M2").Value = -2453124.2
$ws.Range("N2").Value = -2404.3333
$ws.Range("H45").Value = 6041.4
$ws.Range("I45").Value = 8290.223
$ws.Range("J45").Value = 2668.1667
$ws.Range("K45").Value = 8290.223
$ws.Range("L45").Value = 2668.1667
$ws.Range("M45").Value = -7913.223
$ws.Range("N45").Value = -3422.1667
$ws.Range("H116").Value = 1227707.8
$ws.Range("I116").Value = 2453237.2
$ws.Range("J116").Value = 2178.3333
$ws.Range("K116").Value = 2453237.2
$ws.Range("L116").Value = 2178.3333
$ws.Range("M116").Value = -2450943.2
$ws.Range("N116").Value = -6766.3333
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1227707.8
$ws.Range("I3").Value = 2453237.2
$ws.Range("J3").Value = 2178.3333
$ws.Range("K3").Value = 2453237.2
$ws.Range("L3").Value = 2178.3333
$ws.Range("M3").Value = -2453123.2
$ws.Range("N3").Value = -2406.3333
$ws.Range("H86").Value = 2920.9688
$ws.Range("I86").Value = 3049.95
$ws.Range("J86").Value = 2706
$ws.Range("K86").Value = 3049.95
$ws.Range("L86").Value = 2706
$ws.Range("M86").Value = -1926.95
$ws.Range("N86").Value = -4952
$ws.Range("H89").Value = 2920.9688
$ws.Range("I89").Value = 3049.95
$ws.Range("J89").Value = 2706
$ws.Range("K89").Value = 15249.75
$ws.Range("L89").Value = 13530
$ws.Range("M89").Value = -9633.75
$ws.Range("N89").Value = -24762
$ws.Range("H105").Value = 2595.238
$ws.Range("I105").Value = 2158.3635
$ws.Range("K105").Value = 2158.3635
$ws.Range("M105").Value = -411.3634999999999
$ws.Range("H107").Value = 51857.715
$ws.Range("I107").Value = 4195.2
$ws.Range("J107").Value = 171014
$ws.Range("K107").Value = 4195.2
$ws.Range("L107").Value = 171014
$ws.Range("M107").Value = -2275.2
$ws.Range("N107").Value = -174854
$ws.Range("H134").Value = 21742966
$ws.Range("I134").Value = 25003716
$ws.Range("K134").Value = 75011148
$ws.Range("M134").Value = -75008613
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3906
$ws.Range("I31").Value = 4343.4287
$ws.Range("K31").Value = 4343.4287
$ws.Range("M31").Value = -4048.4287
$ws.Range("H32").Value = 19403
$ws.Range("I32").Value = 8940
$ws.Range("K32").Value = 8940
$ws.Range("M32").Value = -8624
$ws.Range("H34").Value = 3906
$ws.Range("I34").Value = 4343.4287
$ws.Range("K34").Value = 4343.4287
$ws.Range("M34").Value = -4141.4287
$ws.Range("H63").Value = 85900
$ws.Range("I63").Value = 85900
$ws.Range("K63").Value = 85900
$ws.Range("M63").Value = -85214
$ws.Range("H66").Value = 85900
$ws.Range("I66").Value = 85900
$ws.Range("K66").Value = 257700
$ws.Range("M66").Value = -254268
$ws.Range("H105").Value = 1572180.8
$ws.Range("I105").Value = 2269705.5
$ws.Range("K105").Value = 2269705.5
$ws.Range("M105").Value = -2267958.5
$ws.Range("H116").Value = 34995.5
$ws.Range("J116").Value = 34995.5
$ws.Range("L116").Value = 34995.5
$ws.Range("N116").Value = -44173.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 441476.88
$ws.Range("I4").Value = 579785.4
$ws.Range("J4").Value = 3500
$ws.Range("K4").Value = 1739356.2
$ws.Range("L4").Value = 10500
$ws.Range("M4").Value = -1739244.2
$ws.Range("N4").Value = -10724
$ws.Range("H5").Value = 500499.5
$ws.Range("I5").Value = 500499.5
$ws.Range("K5").Value = 1501498.5
$ws.Range("M5").Value = -1501386.5
$ws.Range("H26").Value = 8205.308000000001
$ws.Range("I26").Value = 338.125
$ws.Range("J26").Value = 20792.8
$ws.Range("K26").Value = 1014.375
$ws.Range("L26").Value = 62378.39999999999
$ws.Range("M26").Value = -726.375
$ws.Range("N26").Value = -62954.39999999999
$ws.Range("H33").Value = 311.11765
$ws.Range("I33").Value = 267
$ws.Range("J33").Value = 392
$ws.Range("K33").Value = 1602
$ws.Range("L33").Value = 2352
$ws.Range("M33").Value = -1319
$ws.Range("N33").Value = -2918
$ws.Range("H46").Value = 500
$ws.Range("J46").Value = 500
$ws.Range("L46").Value = 1500
$ws.Range("N46").Value = -1682
$ws.Range("H69").Value = 625.75
$ws.Range("I69").Value = 625.75
$ws.Range("K69").Value = 1877.25
$ws.Range("M69").Value = -1066.25
$ws.Range("H72").Value = 625.75
$ws.Range("I72").Value = 625.75
$ws.Range("K72").Value = 5631.75
$ws.Range("M72").Value = -1575.75
$ws.Range("H131").Value = 1814.25
$ws.Range("J131").Value = 2072.4546
$ws.Range("L131").Value = 6217.3638
$ws.Range("N131").Value = -16297.3638
$ws.Range("H135").Value = 500499.5
$ws.Range("I135").Value = 500499.5
$ws.Range("K135").Value = 4504495.5
$ws.Range("M135").Value = -4501960.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 8500
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1555.3478
$ws.Range("J16").Value = 2239.375
$ws.Range("L16").Value = 2239.375
$ws.Range("N16").Value = -2579.375
$ws.Range("H46").Value = 1645.2106
$ws.Range("I46").Value = 1753.6875
$ws.Range("K46").Value = 1753.6875
$ws.Range("M46").Value = -1565.6875
$ws.Range("H93").Value = 2689.7273
$ws.Range("I93").Value = 1174.5
$ws.Range("K93").Value = 1174.5
$ws.Range("M93").Value = 73.5
$ws.Range("H96").Value = 42998.5
$ws.Range("J96").Value = 42998.5
$ws.Range("L96").Value = 42998.5
$ws.Range("N96").Value = -48490.5
$ws.Range("H100").Value = 10506866
$ws.Range("I100").Value = 28513966
$ws.Range("J100").Value = 2724
$ws.Range("K100").Value = 28513966
$ws.Range("L100").Value = 2724
$ws.Range("M100").Value = -28513425
$ws.Range("N100").Value = -3806
$ws.Range("H122").Value = 5381.8423
$ws.Range("I122").Value = 4769.615
$ws.Range("K122").Value = 14308.845
$ws.Range("M122").Value = -11858.845
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 26500
$ws.Range("I70").Value = 26500
$ws.Range("K70").Value = 26500
$ws.Range("M70").Value = -26185
$ws.Range("H73").Value = 26500
$ws.Range("I73").Value = 26500
$ws.Range("K73").Value = 26500
$ws.Range("M73").Value = -25408
$ws.Range("H81").Value = 2221.8333
$ws.Range("I81").Value = 777
$ws.Range("J81").Value = 3666.6667
$ws.Range("K81").Value = 1554
$ws.Range("L81").Value = 7333.3334
$ws.Range("M81").Value = -493
$ws.Range("N81").Value = -9455.3334
$ws.Range("H84").Value = 2221.8333
$ws.Range("I84").Value = 777
$ws.Range("J84").Value = 3666.6667
$ws.Range("K84").Value = 7770
$ws.Range("L84").Value = 36666.667
$ws.Range("M84").Value = -2466
$ws.Range("N84").Value = -47274.667
$ws.Range("H100").Value = 2697.7778
$ws.Range("I100").Value = 2697.7778
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 5395.5556
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -4854.5556
$ws.Range("N100").ClearContents()
